# The document swaps the (cosmetic) "name" metadata carried on the two
# logo pictures that live in the headers/footers:
#   - the Pearson Edexcel logo (alt text / description ends in
#     "PearsonLogo.png") is renamed from "image1.png" -> "image2.png"
#     (it appears twice: in the "first page" footer and the "default"
#     footer)
#   - the BTEC logo (alt text / description "BTec_Logo-Orange") is
#     renamed from "image2.jpg" -> "image1.jpg" (it appears once, in the
#     "first page" header)
# The embedded image data / relationships themselves are untouched -
# only the shape's display name changes.

$d = $word.ActiveDocument

function Rename-LogoShapes($range, $targetDescr, $newName) {
    if ($range -eq $null) { return }
    $shapes = $range.InlineShapes
    if ($shapes -eq $null) { return }
    $count = $shapes.Count
    for ($i = 1; $i -le $count; $i++) {
        $shape = $shapes.Item($i)
        if ($shape.AlternativeText -eq $targetDescr) {
            $shape.Name = $newName
        }
    }
}

foreach ($sec in $d.Sections) {
    # wdHeaderFooterPrimary=1, wdHeaderFooterFirstPage=2, wdHeaderFooterEvenPages=3
    for ($hf = 1; $hf -le 3; $hf++) {
        Rename-LogoShapes $sec.Headers($hf).Range `
            "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" `
            "image2.png"
        Rename-LogoShapes $sec.Footers($hf).Range `
            "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" `
            "image2.png"

        Rename-LogoShapes $sec.Headers($hf).Range "BTec_Logo-Orange" "image1.jpg"
        Rename-LogoShapes $sec.Footers($hf).Range "BTec_Logo-Orange" "image1.jpg"
    }
}
